$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Suit" column (D), which was previously blank, so the
# cards can be organized by suit then rank. The suit letter matches
# the tail of each row's DB ID in column A (…KC.1 -> C, …KD.1 -> D,
# …KH.1 -> H, …KS.1 -> S), repeating for the K, Q and J rank blocks.

$ws.Range("D2").Value  = "C"
$ws.Range("D3").Value  = "D"
$ws.Range("D4").Value  = "H"
$ws.Range("D5").Value  = "S"

$ws.Range("D6").Value  = "C"
$ws.Range("D7").Value  = "D"
$ws.Range("D8").Value  = "H"
$ws.Range("D9").Value  = "S"

$ws.Range("D10").Value = "C"
$ws.Range("D11").Value = "D"
$ws.Range("D12").Value = "H"
$ws.Range("D13").Value = "S"

# Leave the selection where the editor ended up after the last edit.
$ws.Range("D17").Select()
